$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry: 2017-07-22 (serial 42938), 1 hour, with a note about
# switching the FileChooser/config-file UI work over to JavaFX.
# (Values/formulas are written before formatting is copied over, otherwise
# the SUM(C:C) dependency on the newly-touched column C cell doesn't pick
# up live.)
$ws.Range("A30").Value = 42938
$ws.Range("B30").Formula = "=A30"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "Started using the FileChooser class to open configuration files and the text is too small to read. I can't fix this the way I have with the other gui components. Looks like I have to convert the entire application to JavaFX."

# Copy the formatting (number format, alignment, wrap text, etc.) from the
# previously-last data row (29) down onto the new row so the new entry
# matches the rest of the log table, then match its taller wrapped-text row
# height.
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(30).RowHeight = 42.75

# Recalculate so the Weekly/Total Hours Spent formula (F2 = SUM(C:C)) picks
# up the new hour entry.
$excel.Calculate()

# Mirror the author's final cursor position (one row below the new entry).
$ws.Range("D31").Select() | Out-Null
